$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs in-place (rich-text runs keep their own formatting) ---
# A8 shared string: "Volume 32   Number  23" -> "...24" (issue number)
$ws.Range("A8").Characters(21, 2).Text = "24"
# C9 shared string: "Report Covering the Week  6/2/2025  Through  6/8/2025" -> new week range
$ws.Range("C9").Characters(27, 8).Text = "6/9/2025"
$ws.Range("C9").Characters(46, 8).Text = "6/15/2025"

# --- Update weekly crime statistics table (rows 14-31) ---
# Row 14: Murder
$ws.Range("F14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("L14").Value = -40
$ws.Range("N14").Value = -80
# Row 15: Rape
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 16
$ws.Range("K15").Value = -23.809523809523
$ws.Range("L15").Value = -23.809523809523
$ws.Range("M15").Value = 23.076923076923
$ws.Range("N15").Value = -58.974358974359
# Row 16: Robbery
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -70
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 41
$ws.Range("H16").Value = -19.512195121951
$ws.Range("I16").Value = 210
$ws.Range("J16").Value = 244
$ws.Range("K16").Value = -13.934426229508
$ws.Range("L16").Value = -17.647058823529
$ws.Range("M16").Value = -2.777777777777
$ws.Range("N16").Value = -76.846747519294
# Row 17: Fel. Assault
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 29
$ws.Range("E17").Value = -34.482758620689
$ws.Range("F17").Value = 73
$ws.Range("G17").Value = 96
$ws.Range("H17").Value = -23.958333333333
$ws.Range("I17").Value = 432
$ws.Range("J17").Value = 490
$ws.Range("K17").Value = -11.836734693877
$ws.Range("L17").Value = -11.111111111111
$ws.Range("M17").Value = 72.111553784860
$ws.Range("N17").Value = -25.129982668977
# Row 18: Burglary
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 106
$ws.Range("J18").Value = 130
$ws.Range("K18").Value = -18.461538461538
$ws.Range("L18").Value = -55.648535564853
$ws.Range("M18").Value = -9.401709401709
$ws.Range("N18").Value = -88.723404255319
# Row 19: Gr. Larceny
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -9.090909090909
$ws.Range("I19").Value = 346
$ws.Range("J19").Value = 409
$ws.Range("K19").Value = -15.403422982885
$ws.Range("L19").Value = 1.466275659824
$ws.Range("M19").Value = 107.185628742515
$ws.Range("N19").Value = 0
# Row 20: G.L.A.
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -20.689655172413
$ws.Range("I20").Value = 162
$ws.Range("J20").Value = 155
$ws.Range("K20").Value = 4.516129032258
$ws.Range("L20").Value = -28.318584070796
$ws.Range("M20").Value = 88.372093023255
$ws.Range("N20").Value = -75.820895522388
# Row 21: TOTAL
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -22.413793103448
$ws.Range("F21").Value = 218
$ws.Range("G21").Value = 266
$ws.Range("H21").Value = -18.045112781954
$ws.Range("I21").Value = 1278
$ws.Range("J21").Value = 1456
$ws.Range("K21").Value = -12.225274725274
$ws.Range("L21").Value = -19.011406844106
$ws.Range("M21").Value = 48.951048951049
$ws.Range("N21").Value = -63.579367341122
# Row 22: Transit
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = 0
# Row 23: Housing
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -60
$ws.Range("J23").Value = 36
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -29.411764705882
$ws.Range("M23").Value = -7.692307692307
# Row 24: Petit Larceny
$ws.Range("C24").Value = 45
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = 55.172413793103
$ws.Range("F24").Value = 163
$ws.Range("G24").Value = 131
$ws.Range("H24").Value = 24.427480916030
$ws.Range("I24").Value = 891
$ws.Range("J24").Value = 892
$ws.Range("K24").Value = -0.112107623318
$ws.Range("L24").Value = -6.603773584905
$ws.Range("M24").Value = 49.496644295302
# Row 25: Retail Theft
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 116.666666666667
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 29.545454545454
$ws.Range("I25").Value = 325
$ws.Range("J25").Value = 455
$ws.Range("K25").Value = -28.571428571428
$ws.Range("L25").Value = -35
# Row 26: Misd. Assault
$ws.Range("C26").Value = 31
$ws.Range("D26").Value = 34
$ws.Range("E26").Value = -8.823529411764
$ws.Range("F26").Value = 134
$ws.Range("G26").Value = 137
$ws.Range("H26").Value = -2.189781021897
$ws.Range("I26").Value = 572
$ws.Range("J26").Value = 630
$ws.Range("K26").Value = -9.206349206349
$ws.Range("L26").Value = 1.418439716312
$ws.Range("M26").Value = 5.535055350553
# Row 27: UCR Rape*
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("E27").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = -28.125
$ws.Range("L27").Value = -32.352941176470
# Row 28: Other Sex Crimes
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("E28").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("F28").Value = 17
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 183.333333333333
$ws.Range("I28").Value = 69
$ws.Range("K28").Value = 15
$ws.Range("L28").Value = 1.470588235294
# Row 29: Shooting Vic.
$ws.Range("C29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("N29").Value = -87.341772151898
# Row 30: Shooting Inc.
$ws.Range("C30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("N30").Value = -90
# Row 31: Hate Crimes
$ws.Range("D31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)  # restore number format/style for type change
$ws.Range("E31").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)  # restore number format/style for type change

$excel.CutCopyMode = $false
